# Weekly refresh of Fruta/Hortalizas price data (Kiwi, Agricola del Norte S.A. de Arica).
# Columns D (Fecha), K (Variedad), L (Calidad), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), Q (Unidad de comercializacion),
# R (Origen), S (Precio $/Kg) and T (Kg / unidad) are refreshed per row 2-29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    [pscustomobject]@{ Row = 2; Fecha = 45134; Variedad = 'Hayward'; Calidad = 'Especial'; Volumen = 350; PrecioMin = 21000; PrecioMax = 22000; PrecioProm = 21429; Unidad = '$/bandeja 18 kilos'; Origen = 'Región de O''Higgins'; PrecioKg = 1190; KgUnidad = 18 },
    [pscustomobject]@{ Row = 3; Fecha = 45086; Variedad = 'Hayward'; Calidad = 'Especial'; Volumen = 250; PrecioMin = 25000; PrecioMax = 26000; PrecioProm = 25500; Unidad = '$/bandeja 18 kilos'; Origen = 'Región de O''Higgins'; PrecioKg = 1417; KgUnidad = 18 },
    [pscustomobject]@{ Row = 4; Fecha = 45086; Variedad = 'Hayward'; Calidad = 'Primera'; Volumen = 250; PrecioMin = 20000; PrecioMax = 21000; PrecioProm = 20500; Unidad = '$/bandeja 18 kilos'; Origen = 'Región de O''Higgins'; PrecioKg = 1139; KgUnidad = 18 },
    [pscustomobject]@{ Row = 5; Fecha = 44991; Variedad = 'Hayward'; Calidad = 'Primera'; Volumen = 250; PrecioMin = 24000; PrecioMax = 25000; PrecioProm = 24500; Unidad = '$/bandeja 18 kilos'; Origen = 'Región de O''Higgins'; PrecioKg = 1361; KgUnidad = 18 },
    [pscustomobject]@{ Row = 6; Fecha = 44789; Variedad = 'Hayward'; Calidad = 'Segunda'; Volumen = 250; PrecioMin = 19000; PrecioMax = 20000; PrecioProm = 19500; Unidad = '$/bandeja 18 kilos'; Origen = 'Región de O''Higgins'; PrecioKg = 1083; KgUnidad = 18 },
    [pscustomobject]@{ Row = 7; Fecha = 44656; Variedad = 'Hayward'; Calidad = 'Primera'; Volumen = 270; PrecioMin = 19000; PrecioMax = 20000; PrecioProm = 19500; Unidad = '$/bandeja 18 kilos'; Origen = 'Región de O''Higgins'; PrecioKg = 1083; KgUnidad = 18 },
    [pscustomobject]@{ Row = 8; Fecha = 44784; Variedad = 'Hayward'; Calidad = 'Primera'; Volumen = 300; PrecioMin = 19000; PrecioMax = 20000; PrecioProm = 19500; Unidad = '$/bandeja 18 kilos'; Origen = 'Región de O''Higgins'; PrecioKg = 1083; KgUnidad = 18 },
    [pscustomobject]@{ Row = 9; Fecha = 45002; Variedad = 'Hayward'; Calidad = 'Segunda'; Volumen = 300; PrecioMin = 24000; PrecioMax = 25000; PrecioProm = 24500; Unidad = '$/bandeja 18 kilos'; Origen = 'Región de O''Higgins'; PrecioKg = 1361; KgUnidad = 18 },
    [pscustomobject]@{ Row = 10; Fecha = 44323; Variedad = 'Hayward'; Calidad = 'Primera'; Volumen = 270; PrecioMin = 21000; PrecioMax = 22000; PrecioProm = 21500; Unidad = '$/bandeja 18 kilos'; Origen = 'Región de O''Higgins'; PrecioKg = 1194; KgUnidad = 18 },
    [pscustomobject]@{ Row = 11; Fecha = 44307; Variedad = 'Hayward'; Calidad = 'Primera'; Volumen = 250; PrecioMin = 19000; PrecioMax = 20000; PrecioProm = 19500; Unidad = '$/bandeja 18 kilos'; Origen = 'Región de O''Higgins'; PrecioKg = 1083; KgUnidad = 18 },
    [pscustomobject]@{ Row = 12; Fecha = 45107; Variedad = 'Hayward'; Calidad = 'Primera'; Volumen = 320; PrecioMin = 20000; PrecioMax = 21000; PrecioProm = 20500; Unidad = '$/bandeja 18 kilos'; Origen = 'Región de O''Higgins'; PrecioKg = 1139; KgUnidad = 18 },
    [pscustomobject]@{ Row = 13; Fecha = 44629; Variedad = 'Hayward'; Calidad = 'Segunda'; Volumen = 300; PrecioMin = 17000; PrecioMax = 18000; PrecioProm = 17500; Unidad = '$/bandeja 18 kilos'; Origen = 'Región de O''Higgins'; PrecioKg = 972; KgUnidad = 18 },
    [pscustomobject]@{ Row = 14; Fecha = 44263; Variedad = 'Hayward'; Calidad = 'Primera'; Volumen = 250; PrecioMin = 21000; PrecioMax = 22000; PrecioProm = 21500; Unidad = '$/caja 18 kilos'; Origen = 'Región de O''Higgins'; PrecioKg = 1194; KgUnidad = 18 },
    [pscustomobject]@{ Row = 15; Fecha = 44602; Variedad = 'Hayward'; Calidad = 'Primera'; Volumen = 270; PrecioMin = 20000; PrecioMax = 21000; PrecioProm = 20500; Unidad = '$/bandeja 18 kilos'; Origen = 'Región de O''Higgins'; PrecioKg = 1139; KgUnidad = 18 },
    [pscustomobject]@{ Row = 16; Fecha = 45034; Variedad = 'Hayward'; Calidad = 'Primera'; Volumen = 250; PrecioMin = 25000; PrecioMax = 26000; PrecioProm = 25600; Unidad = '$/bandeja 18 kilos'; Origen = 'Región de O''Higgins'; PrecioKg = 1422; KgUnidad = 18 },
    [pscustomobject]@{ Row = 17; Fecha = 45127; Variedad = 'Hayward'; Calidad = 'Primera'; Volumen = 200; PrecioMin = 21000; PrecioMax = 22000; PrecioProm = 21500; Unidad = '$/bandeja 18 kilos'; Origen = 'Región de O''Higgins'; PrecioKg = 1194; KgUnidad = 18 },
    [pscustomobject]@{ Row = 18; Fecha = 44487; Variedad = 'Hayward'; Calidad = 'Primera'; Volumen = 300; PrecioMin = 14000; PrecioMax = 15000; PrecioProm = 14500; Unidad = '$/bandeja 10 kilos'; Origen = 'Región de O''Higgins'; PrecioKg = 1450; KgUnidad = 10 },
    [pscustomobject]@{ Row = 19; Fecha = 44706; Variedad = 'Hayward'; Calidad = 'Primera'; Volumen = 400; PrecioMin = 9000; PrecioMax = 10000; PrecioProm = 9500; Unidad = '$/bandeja 10 kilos'; Origen = 'Región de O''Higgins'; PrecioKg = 950; KgUnidad = 10 },
    [pscustomobject]@{ Row = 20; Fecha = 44418; Variedad = 'Hayward'; Calidad = 'Primera'; Volumen = 240; PrecioMin = 10000; PrecioMax = 11000; PrecioProm = 10500; Unidad = '$/bandeja 10 kilos'; Origen = 'Región de O''Higgins'; PrecioKg = 1050; KgUnidad = 10 },
    [pscustomobject]@{ Row = 21; Fecha = 44673; Variedad = 'Hayward'; Calidad = 'Especial'; Volumen = 400; PrecioMin = 14000; PrecioMax = 15000; PrecioProm = 14500; Unidad = '$/bandeja 10 kilos'; Origen = 'Región de O''Higgins'; PrecioKg = 1450; KgUnidad = 10 },
    [pscustomobject]@{ Row = 22; Fecha = 44614; Variedad = 'Hayward'; Calidad = 'Primera'; Volumen = 250; PrecioMin = 20000; PrecioMax = 21000; PrecioProm = 20500; Unidad = '$/bandeja 18 kilos'; Origen = 'Región de O''Higgins'; PrecioKg = 1139; KgUnidad = 18 },
    [pscustomobject]@{ Row = 23; Fecha = 45069; Variedad = 'Sin especificar'; Calidad = 'Primera'; Volumen = 370; PrecioMin = 19000; PrecioMax = 20000; PrecioProm = 19486; Unidad = '$/bandeja 18 kilos'; Origen = 'Región Metropolitana'; PrecioKg = 1083; KgUnidad = 18 },
    [pscustomobject]@{ Row = 24; Fecha = 44489; Variedad = 'Hayward'; Calidad = 'Primera'; Volumen = 300; PrecioMin = 26000; PrecioMax = 27000; PrecioProm = 26500; Unidad = '$/bandeja 18 kilos'; Origen = 'Región de O''Higgins'; PrecioKg = 1472; KgUnidad = 18 },
    [pscustomobject]@{ Row = 25; Fecha = 44491; Variedad = 'Hayward'; Calidad = 'Primera'; Volumen = 300; PrecioMin = 14000; PrecioMax = 15000; PrecioProm = 14500; Unidad = '$/bandeja 10 kilos'; Origen = 'Región de O''Higgins'; PrecioKg = 1450; KgUnidad = 10 },
    [pscustomobject]@{ Row = 26; Fecha = 44291; Variedad = 'Hayward'; Calidad = 'Primera'; Volumen = 200; PrecioMin = 17000; PrecioMax = 18000; PrecioProm = 17500; Unidad = '$/bandeja 18 kilos'; Origen = 'Región de O''Higgins'; PrecioKg = 972; KgUnidad = 18 },
    [pscustomobject]@{ Row = 27; Fecha = 44819; Variedad = 'Hayward'; Calidad = 'Primera'; Volumen = 300; PrecioMin = 17000; PrecioMax = 18000; PrecioProm = 17500; Unidad = '$/bandeja 10 kilos'; Origen = 'Región de O''Higgins'; PrecioKg = 1750; KgUnidad = 10 },
    [pscustomobject]@{ Row = 28; Fecha = 45043; Variedad = 'Hayward'; Calidad = 'Segunda'; Volumen = 300; PrecioMin = 21000; PrecioMax = 22000; PrecioProm = 21500; Unidad = '$/bandeja 18 kilos'; Origen = 'Región de O''Higgins'; PrecioKg = 1194; KgUnidad = 18 },
    [pscustomobject]@{ Row = 29; Fecha = 44616; Variedad = 'Hayward'; Calidad = 'Segunda'; Volumen = 300; PrecioMin = 16000; PrecioMax = 17000; PrecioProm = 16500; Unidad = '$/caja 18 kilos granel'; Origen = 'Región de O''Higgins'; PrecioKg = 917; KgUnidad = 18 }
)

foreach ($row in $rows) {
    $ws.Cells.Item($row.Row, 4).Value = $row.Fecha
    $ws.Cells.Item($row.Row, 11).Value = $row.Variedad
    $ws.Cells.Item($row.Row, 12).Value = $row.Calidad
    $ws.Cells.Item($row.Row, 13).Value = $row.Volumen
    $ws.Cells.Item($row.Row, 14).Value = $row.PrecioMin
    $ws.Cells.Item($row.Row, 15).Value = $row.PrecioMax
    $ws.Cells.Item($row.Row, 16).Value = $row.PrecioProm
    $ws.Cells.Item($row.Row, 17).Value = $row.Unidad
    $ws.Cells.Item($row.Row, 18).Value = $row.Origen
    $ws.Cells.Item($row.Row, 19).Value = $row.PrecioKg
    $ws.Cells.Item($row.Row, 20).Value = $row.KgUnidad
}
